$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 3246931
$ws.Range("I12").Value = 6493594.5
$ws.Range("J12").Value = 267.57144
$ws.Range("K12").Value = 6493594.5
$ws.Range("L12").Value = 267.57144
$ws.Range("M12").Value = -6493424.5
$ws.Range("N12").Value = -607.5714399999999
# Row 21
$ws.Range("H21").Value = 5000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents() | Out-Null
# Row 23
$ws.Range("H23").Value = 5000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents() | Out-Null
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents() | Out-Null
# Row 32
$ws.Range("H32").Value = 55557310
$ws.Range("I32").Value = 1509
$ws.Range("K32").Value = 1509
$ws.Range("M32").Value = -1183
# Row 70
$ws.Range("H70").Value = 2758.5715
$ws.Range("J70").Value = 3439.8
$ws.Range("L70").Value = 10319.4
$ws.Range("N70").Value = -10859.4
# Row 73
$ws.Range("H73").Value = 2758.5715
$ws.Range("J73").Value = 3439.8
$ws.Range("L73").Value = 10319.4
$ws.Range("N73").Value = -12191.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 40
$ws.Range("H40").Value = 20000
$ws.Range("J40").Value = 20000
$ws.Range("L40").Value = 20000
$ws.Range("N40").Value = -20352
# Row 42
$ws.Range("H42").Value = 12345
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 12345
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 12345
$ws.Range("M42").ClearContents() | Out-Null
$ws.Range("N42").Value = -13317
# Row 45
$ws.Range("H45").Value = 144985.14
$ws.Range("I45").Value = 168649.33
$ws.Range("K45").Value = 168649.33
$ws.Range("M45").Value = -168272.33
# Row 61
$ws.Range("H61").Value = 849817.1
$ws.Range("I61").Value = 24983.959
$ws.Range("J61").Value = 2607070.2
$ws.Range("K61").Value = 24983.959
$ws.Range("L61").Value = 2607070.2
$ws.Range("M61").Value = -24771.959
$ws.Range("N61").Value = -2607494.2
# Row 136
$ws.Range("H136").Value = 849817.1
$ws.Range("I136").Value = 24983.959
$ws.Range("J136").Value = 2607070.2
$ws.Range("K136").Value = 74951.87699999999
$ws.Range("L136").Value = 7821210.600000001
$ws.Range("M136").Value = -72401.87699999999
$ws.Range("N136").Value = -7826310.600000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 24
$ws.Range("H24").Value = 7875
$ws.Range("J24").Value = 10166.667
$ws.Range("L24").Value = 10166.667
$ws.Range("N24").Value = -10636.667
# Row 105
$ws.Range("H105").Value = 10417.871
$ws.Range("I105").Value = 6663.32
$ws.Range("K105").Value = 6663.32
$ws.Range("M105").Value = -4916.32
# Row 107
$ws.Range("H107").Value = 7357.5454
$ws.Range("I107").Value = 8057.724
$ws.Range("K107").Value = 8057.724
$ws.Range("M107").Value = -6137.724

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 28
$ws.Range("H28").Value = 31881
$ws.Range("J28").Value = 31881
$ws.Range("L28").Value = 31881
$ws.Range("N28").Value = -32371
# Row 50
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents() | Out-Null
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents() | Out-Null
# Row 58
$ws.Range("H58").Value = 1474.7941
$ws.Range("I58").Value = 1318.0555
$ws.Range("J58").Value = 1651.125
$ws.Range("K58").Value = 1318.0555
$ws.Range("L58").Value = 1651.125
$ws.Range("M58").Value = -1115.0555
$ws.Range("N58").Value = -2057.125
# Row 59
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents() | Out-Null
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents() | Out-Null
# Row 68
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents() | Out-Null
# Row 71
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents() | Out-Null
# Row 80
$ws.Range("H80").Value = 77000
$ws.Range("J80").Value = 77000
$ws.Range("L80").Value = 77000
$ws.Range("N80").Value = -79246
# Row 83
$ws.Range("H83").Value = 77000
$ws.Range("J83").Value = 77000
$ws.Range("L83").Value = 231000
$ws.Range("N83").Value = -242232
# Row 133
$ws.Range("H133").Value = 73123
$ws.Range("J133").Value = 73123
$ws.Range("L133").Value = 73123
$ws.Range("N133").Value = -78183
# Row 134
$ws.Range("H134").Value = 2535.72
$ws.Range("I134").Value = 2313.6
$ws.Range("J134").Value = 2868.9
$ws.Range("K134").Value = 6940.799999999999
$ws.Range("L134").Value = 8606.700000000001
$ws.Range("M134").Value = -4405.799999999999
$ws.Range("N134").Value = -13676.7
# Row 136
$ws.Range("H136").Value = 1474.7941
$ws.Range("I136").Value = 1318.0555
$ws.Range("J136").Value = 1651.125
$ws.Range("K136").Value = 3954.1665
$ws.Range("L136").Value = 4953.375
$ws.Range("M136").Value = -1404.1665
$ws.Range("N136").Value = -10053.375
# Row 137
$ws.Range("H137").Value = 80000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents() | Out-Null

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 2357908.5
$ws.Range("I4").Value = 4231519
$ws.Range("J4").Value = 734112.9
$ws.Range("K4").Value = 12694557
$ws.Range("L4").Value = 2202338.7
$ws.Range("M4").Value = -12694445
$ws.Range("N4").Value = -2202562.7
# Row 44
$ws.Range("H44").Value = 6666.4
$ws.Range("J44").Value = 6666.4
$ws.Range("L44").Value = 19999.2
$ws.Range("N44").Value = -20795.2
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents() | Out-Null
$ws.Range("N57").ClearContents() | Out-Null
# Row 103
$ws.Range("H103").Value = 669952.9399999999
$ws.Range("I103").Value = 1666882.1
$ws.Range("J103").Value = 5333.4443
$ws.Range("K103").Value = 5000646.300000001
$ws.Range("L103").Value = 16000.3329
$ws.Range("M103").Value = -4999767.300000001
$ws.Range("N103").Value = -17758.3329
# Row 114
$ws.Range("H114").Value = 1416.2106
$ws.Range("I114").Value = 1401.5
$ws.Range("J114").Value = 1423
$ws.Range("K114").Value = 4204.5
$ws.Range("L114").Value = 4269
$ws.Range("M114").Value = -950.5
$ws.Range("N114").Value = -10777
# Row 121
$ws.Range("H121").Value = 2004.6666
$ws.Range("I121").Value = 405.7
$ws.Range("K121").Value = 1217.1
$ws.Range("M121").Value = 92.90000000000009
# Row 122
$ws.Range("H122").Value = 8334534.5
$ws.Range("J122").Value = 2424.6667
$ws.Range("L122").Value = 21822.0003
$ws.Range("N122").Value = -26722.0003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents() | Out-Null
# Row 126
$ws.Range("H126").Value = 3407.5
$ws.Range("I126").Value = 2200
$ws.Range("K126").Value = 6600
$ws.Range("M126").Value = -4130
# Row 132
$ws.Range("H132").Value = 1074288.2
$ws.Range("J132").Value = 1509730.8
$ws.Range("L132").Value = 4529192.4
$ws.Range("N132").Value = -4534252.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 98
$ws.Range("H98").Value = 66403
$ws.Range("J98").Value = 66403
$ws.Range("L98").Value = 66403
$ws.Range("N98").Value = -72393
# Row 112
$ws.Range("H112").Value = 148777
$ws.Range("J112").Value = 148777
$ws.Range("L112").Value = 148777
$ws.Range("N112").Value = -151731
# Row 132
$ws.Range("H132").Value = 3455.8823
$ws.Range("I132").Value = 3610.6
$ws.Range("J132").Value = 3234.8572
$ws.Range("K132").Value = 10831.8
$ws.Range("L132").Value = 9704.571599999999
$ws.Range("M132").Value = -8301.799999999999
$ws.Range("N132").Value = -14764.5716
# Row 134
$ws.Range("H134").Value = 69982
$ws.Range("J134").Value = 69982
$ws.Range("L134").Value = 69982
$ws.Range("N134").Value = -80122
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents() | Out-Null

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 34239.4
$ws.Range("J70").Value = 36533
$ws.Range("L70").Value = 36533
$ws.Range("N70").Value = -37163
# Row 73
$ws.Range("H73").Value = 34239.4
$ws.Range("J73").Value = 36533
$ws.Range("L73").Value = 36533
$ws.Range("N73").Value = -38717
# Row 97
$ws.Range("H97").Value = 36871
$ws.Range("J97").Value = 36871
$ws.Range("L97").Value = 36871
$ws.Range("N97").Value = -38853
# Row 122
$ws.Range("H122").Value = 1134.7561
$ws.Range("I122").Value = 1109.8064
$ws.Range("K122").Value = 3329.4192
$ws.Range("M122").Value = -879.4191999999998
# Row 132
$ws.Range("H132").Value = 2049.9167
$ws.Range("I132").Value = 1529.2273
$ws.Range("K132").Value = 4587.6819
$ws.Range("M132").Value = -2057.6819
